$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 4,23
$data[0,0] = 0.013555787278415
$data[0,1] = 0.00208550573514077
$data[0,2] = 0.86444212721585
$data[0,3] = 0.00834202294056309
$data[0,4] = 0
$data[0,5] = 0.0312825860271116
$data[0,6] = 0.899895724713243
$data[0,7] = 0.00208550573514077
$data[0,8] = 0.0125130344108446
$data[0,9] = 0
$data[0,10] = 0.00521376433785193
$data[0,11] = 0
$data[0,12] = 0.0114702815432742
$data[0,13] = 0.00417101147028154
$data[0,14] = 0.843587069864442
$data[0,15] = 0.00104275286757039
$data[0,16] = 0.0166840458811262
$data[0,17] = 0.0281543274244004
$data[0,18] = 0.00938477580813347
$data[0,19] = 0.00417101147028154
$data[0,20] = 0.0114702815432742
$data[0,21] = 0.0177267987486966
$data[0,22] = 0.222106360792492
$data[1,0] = 0.962460896767466
$data[1,1] = 0.969760166840459
$data[1,2] = 0.116788321167883
$data[1,3] = 0.905109489051095
$data[1,4] = 0.00521376433785193
$data[1,5] = 0
$data[1,6] = 0
$data[1,7] = 0.896767466110532
$data[1,8] = 0.949947862356621
$data[1,9] = 0.0114702815432742
$data[1,10] = 0.0104275286757039
$data[1,11] = 0.0166840458811262
$data[1,12] = 0.0145985401459854
$data[1,13] = 0.0156412930135558
$data[1,14] = 0.02711157455683
$data[1,15] = 0.0218978102189781
$data[1,16] = 0.00625651720542232
$data[1,17] = 0.00208550573514077
$data[1,18] = 0.0250260688216893
$data[1,19] = 0.0156412930135558
$data[1,20] = 0.94681960375391
$data[1,21] = 0.0114702815432742
$data[1,22] = 0.0239833159541189
$data[2,0] = 0.0156412930135558
$data[2,1] = 0.0198123044838373
$data[2,2] = 0.0125130344108446
$data[2,3] = 0.0166840458811262
$data[2,4] = 0.00208550573514077
$data[2,5] = 0.965589155370177
$data[2,6] = 0.0114702815432742
$data[2,7] = 0.0959332638164755
$data[2,8] = 0
$data[2,9] = 0
$data[2,10] = 0.0208550573514077
$data[2,11] = 0.0333680917622523
$data[2,12] = 0.97393117831074
$data[2,13] = 0.966631908237748
$data[2,14] = 0.122002085505735
$data[2,15] = 0.0198123044838373
$data[2,16] = 0.968717413972888
$data[2,17] = 0.953076120959333
$data[2,18] = 0.0187695516162669
$data[2,19] = 0.962460896767466
$data[2,20] = 0.0291970802919708
$data[2,21] = 0.0187695516162669
$data[2,22] = 0.0187695516162669
$data[3,0] = 0.0072992700729927
$data[3,1] = 0.00834202294056309
$data[3,2] = 0.00625651720542232
$data[3,3] = 0.0698644421272158
$data[3,4] = 0.992700729927007
$data[3,5] = 0.00312825860271116
$data[3,6] = 0.0886339937434828
$data[3,7] = 0.00521376433785193
$data[3,8] = 0.0375391032325339
$data[3,9] = 0.988529718456726
$data[3,10] = 0.963503649635037
$data[3,11] = 0.949947862356621
$data[3,12] = 0
$data[3,13] = 0.013555787278415
$data[3,14] = 0.0072992700729927
$data[3,15] = 0.957247132429614
$data[3,16] = 0.00834202294056309
$data[3,17] = 0.0166840458811262
$data[3,18] = 0.94681960375391
$data[3,19] = 0.0177267987486966
$data[3,20] = 0.0125130344108446
$data[3,21] = 0.952033368091762
$data[3,22] = 0.735140771637122

$ws.Range("B2:X5").Value = $data
